# ------------------------------------------------------------------
# "Updated cp and ic" - adds a new CT-IP deposit row, re-threads the
# CT-IP totals, and fills in the ICICI sheet with real cp/ic data
# rows (previously mostly blank placeholders), plus refreshes the
# active-sheet/selection bookkeeping.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ====================================================================
# Sheet "CT-IP" (sheet1): insert a new deposit row (old row 16/18/20
# shift down to 17/19/21), fill it in, and re-point the two summary
# formulas at the new last data row.
# ====================================================================
$ws1 = $wb.Worksheets.Item(1)

$ws1.Rows.Item(15).Insert()

$ws1.Cells.Item(15,1).Value2 = "CT-IP"
$ws1.Cells.Item(15,2).Value2 = 30000
$ws1.Cells.Item(15,3).Value2 = 41975
$ws1.Cells.Item(15,4).Value2 = 7.75
$ws1.Cells.Item(15,5).Value2 = 61
$ws1.Cells.Item(15,6).Value2 = 395
$ws1.Cells.Item(15,7).Value2 = 42037
$ws1.Cells.Item(15,8).Value2 = 30395
$ws1.Cells.Item(15,9).Value2 = "0X0XX8501300"
$ws1.Cells.Item(15,10).Formula = "=H15-B15"

# old row16 (now row17) SUM needs to cover the new row
$ws1.Cells.Item(17,10).Formula = "=SUM(J3:J15)"
# old row20 (now row21) grand total formula references the renumbered rows
$ws1.Cells.Item(21,10).Formula = "=J17+J19*2"

# ====================================================================
# Sheet "CT-NL" (sheet2): just a selection/active-tab bookkeeping
# change - no data changed.
# ====================================================================
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F30").Select()

# ====================================================================
# Sheet "ICICI" (sheet3): populate cp/ic rows 3-5, restructure the
# summary block (rows 6-17).
# ====================================================================
$ws3 = $wb.Worksheets.Item(3)

# Blank style template (all s="4") used to stamp freshly-used rows.
$ws3.Range("A2:M2").Copy()
$ws3.Range("A6:M6").PasteSpecial(-4122)
$ws3.Range("A7:M7").PasteSpecial(-4122)
$ws3.Range("A8:M8").PasteSpecial(-4122)
$ws3.Range("A9:M9").PasteSpecial(-4122)
$ws3.Range("A10:M10").PasteSpecial(-4122)

$ws3.Cells.Item(2,10).Copy()
$ws3.Cells.Item(12,10).PasteSpecial(-4122)

# Pull in date (s=5) / computed-amount (s=8) formats from sheet "CT-IP"
# so the new cells share the existing named styles instead of minting
# new ones.
$ws1.Cells.Item(3,3).Copy()   # date format (s=5)
$ws3.Range("C3:C5").PasteSpecial(-4122)
$ws3.Range("G3:G5").PasteSpecial(-4122)

$ws1.Cells.Item(3,8).Copy()   # computed-amount format (s=8)
$ws3.Range("F3:F5").PasteSpecial(-4122)
$ws3.Range("H3:H5").PasteSpecial(-4122)
$ws3.Cells.Item(3,10).PasteSpecial(-4122)
$ws3.Range("J4:J5").PasteSpecial(-4122)
$ws3.Cells.Item(7,10).PasteSpecial(-4122)

# --- Row 3 ---
$ws3.Cells.Item(3,3).Value2 = 41925
$ws3.Cells.Item(3,6).Formula = "= ( ( (B3*D3/100) / 365 ) * E3 )"
$ws3.Cells.Item(3,7).Value2 = 42313
$ws3.Cells.Item(3,8).Formula = "=B3+F3"
$ws3.Cells.Item(3,9).Value2 = 4710071983
$ws3.Cells.Item(3,10).Value2 = 2500

# --- Row 4 ---
$ws3.Cells.Item(4,1).Value2 = "IC"
$ws3.Cells.Item(4,2).Value2 = 25000
$ws3.Cells.Item(4,3).Value2 = 41947
$ws3.Cells.Item(4,4).Value2 = 7.75
$ws3.Cells.Item(4,5).Value2 = 61
$ws3.Cells.Item(4,6).Formula = "= ( ( (B4*D4/100) / 365 ) * E4 )"
$ws3.Cells.Item(4,7).Value2 = 42009
$ws3.Cells.Item(4,8).Formula = "=B4+F4"
$ws3.Cells.Item(4,9).Value2 = 4710072323
$ws3.Cells.Item(4,10).Formula = "=F4"
$ws3.Cells.Item(4,13).Value2 = "Monthly"

# --- Row 5 ---
$ws3.Cells.Item(5,1).Value2 = "IC"
$ws3.Cells.Item(5,2).Value2 = 25000
$ws3.Cells.Item(5,3).Value2 = 41975
$ws3.Cells.Item(5,4).Value2 = 8.75
$ws3.Cells.Item(5,5).Value2 = 390
$ws3.Cells.Item(5,6).Formula = "= ( ( (B5*D5/100) / 365 ) * E5 )"
$ws3.Cells.Item(5,8).Formula = "=B5+F5"
$ws3.Cells.Item(5,9).Value2 = 4713004066
$ws3.Cells.Item(5,10).Formula = "=F5"
$ws3.Cells.Item(5,13).Value2 = "Monthly"

# --- Row 7: sub-total of the two cp/ic deposits above ---
$ws3.Cells.Item(7,10).Formula = "=J3+J5"

# --- Row 9 / 10: the two interest figures that used to live on row 5 ---
$ws3.Cells.Item(9,9).Value2 = "INT 2013 QT-1"
$ws3.Cells.Item(9,10).Value2 = 744
$ws3.Cells.Item(10,9).Value2 = "INT 2013 QT-2"
$ws3.Cells.Item(10,10).Value2 = 372

# --- Row 12: sum of the two interest figures ---
$ws3.Cells.Item(12,10).Formula = "=J9+J10"

# --- Row 14: grand total (style s=1, bold/centred) ---
$ws3.Cells.Item(1,1).Copy()    # any s="1" cell on this sheet, format only
$ws3.Cells.Item(14,10).PasteSpecial(-4122)
$ws3.Cells.Item(14,10).Formula = "= J7 + J12*2"

# --- Row 17: "Uptd" marker with its date ---
$ws3.Cells.Item(1,1).Copy()
$ws3.Cells.Item(17,1).PasteSpecial(-4122)
$ws3.Cells.Item(17,1).Value2 = "Uptd"
$ws1.Cells.Item(3,3).Copy()    # date format (s=5) again, for B17
$ws3.Cells.Item(17,2).PasteSpecial(-4122)
$ws3.Cells.Item(17,2).Value2 = 41975

$ws3.Range("F18").Select()

# ====================================================================
# Finish on "CT-IP" so it ends up the active tab / selected cell.
# ====================================================================
$ws1.Range("J15").Select()
